# Update countries & provincias Spain
# Re-sort the province rows that moved position (new case counts shuffled
# their rank in the "Casos totales" table) and refresh the "last updated"
# timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = "Ciudad Real"
$ws.Cells.Item(11, 2).Value = 1422
$ws.Cells.Item(11, 3).Value = 153
$ws.Cells.Item(11, 4).Value = 1273
$ws.Cells.Item(11, 5).Value = 89

$ws.Cells.Item(12, 1).Value = "La Rioja"
$ws.Cells.Item(12, 2).Value = 1236
$ws.Cells.Item(12, 3).Value = 62
$ws.Cells.Item(12, 4).Value = 1119
$ws.Cells.Item(12, 5).Value = 55

$ws.Cells.Item(13, 1).Value = "Toledo"
$ws.Cells.Item(13, 2).Value = 1112
$ws.Cells.Item(13, 3).Value = 153
$ws.Cells.Item(13, 4).Value = 972
$ws.Cells.Item(13, 5).Value = 90

$ws.Cells.Item(14, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(14, 2).Value = 1093
$ws.Cells.Item(14, 3).Value = 19
$ws.Cells.Item(14, 4).Value = 983
$ws.Cells.Item(14, 5).Value = 91

$ws.Cells.Item(15, 1).Value = "A Coruña"
$ws.Cells.Item(15, 2).Value = 1049
$ws.Cells.Item(15, 3).Value = 67
$ws.Cells.Item(15, 4).Value = 1000
$ws.Cells.Item(15, 5).Value = 27

$ws.Cells.Item(16, 1).Value = "Zaragoza"
$ws.Cells.Item(16, 2).Value = 1045
$ws.Cells.Item(16, 3).Value = 68
$ws.Cells.Item(16, 4).Value = 928
$ws.Cells.Item(16, 5).Value = 49

$ws.Cells.Item(17, 1).Value = "Malaga"
$ws.Cells.Item(17, 2).Value = 1006
$ws.Cells.Item(17, 3).Value = 61
$ws.Cells.Item(17, 4).Value = 899
$ws.Cells.Item(17, 5).Value = 46

$ws.Cells.Item(25, 1).Value = "Albacete"
$ws.Cells.Item(25, 2).Value = 780
$ws.Cells.Item(25, 3).Value = 153
$ws.Cells.Item(25, 4).Value = 667
$ws.Cells.Item(25, 5).Value = 83

$ws.Cells.Item(26, 1).Value = "Granada"
$ws.Cells.Item(26, 2).Value = 711
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 676
$ws.Cells.Item(26, 5).Value = 34

$ws.Cells.Item(27, 1).Value = "Sevilla"
$ws.Cells.Item(27, 2).Value = 708
$ws.Cells.Item(27, 3).Value = 8
$ws.Cells.Item(27, 4).Value = 675
$ws.Cells.Item(27, 5).Value = 25

$ws.Cells.Item(28, 1).Value = "Murcia"
$ws.Cells.Item(28, 2).Value = 687
$ws.Cells.Item(28, 3).Value = 12
$ws.Cells.Item(28, 4).Value = 660
$ws.Cells.Item(28, 5).Value = 15

$ws.Cells.Item(34, 1).Value = "Guadalajara"
$ws.Cells.Item(34, 2).Value = 440
$ws.Cells.Item(34, 3).Value = 153
$ws.Cells.Item(34, 4).Value = 362
$ws.Cells.Item(34, 5).Value = 75

$ws.Cells.Item(45, 1).Value = "Cuenca"
$ws.Cells.Item(45, 2).Value = 180
$ws.Cells.Item(45, 3).Value = 153
$ws.Cells.Item(45, 4).Value = 130
$ws.Cells.Item(45, 5).Value = 40

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 14:42"
